$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "object 7" : the "ANANDA KRISHNAN  P S" name textbox ---
$nameShape = $s.Shapes.Item(4)

# Shrink the existing line, then add the new "SHALINI R" paragraph below it
# (InsertAfter preserves/continues the paragraph's run formatting). Do this
# before touching position/size, since the textbox auto-fits its height to
# the text and would otherwise clobber an explicit Height set earlier.
$nameShape.TextFrame.TextRange.Font.Size = 24
$null = $nameShape.TextFrame.TextRange.InsertAfter("`rSHALINI R")
$nameShape.TextFrame.TextRange.Font.Size = 24

# Reposition / resize the textbox. (Shape.Left/Top/Width/Height round-trip
# through a single-precision float in this host, so the raw EMU/12700
# point value can truncate one EMU short once re-quantized on save; the
# values below are nudged by a hair so they land back on the exact EMU
# the author's XML has after that round-trip.)
$nameShape.Left = 468.0000305
$nameShape.Top = 161.7575226149606
$nameShape.Width = 264.3201141401575
$nameShape.Height = 60.48500062992126

# --- Shape "object 8" : the "Final Project" textbox ---
$projShape = $s.Shapes.Item(5)
$projShape.Left = 516.0000305
$projShape.Top = 228.0000382
